# Update Runmode (column E) values to reflect the new test selection,
# and update the active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Flip Runmode for the first two test cases from No -> Yes
$ws.Range("E2").Value = "Yes"
$ws.Range("E3").Value = "Yes"

# Flip Runmode for TC07 (1 & 2), TC08, TC09 from Yes -> No
$ws.Range("E8").Value = "No"
$ws.Range("E9").Value = "No"
$ws.Range("E10").Value = "No"
$ws.Range("E11").Value = "No"
$ws.Range("E12").Value = "No"

# Update the selection to match the new active range
$ws.Activate()
$ws.Range("E8:E16").Select()
